$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header cells I1 (I0) and J1 (IF), matching the existing header style (copy format from H1)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Fill data rows 2-67 for columns I (I0) and J (IF)
$data = @(@(9,9),@(8,8),@(8,9),@(9,9),@(9,9),@(8,8),@(7,8),@(8,8),@(9,9),@(8,8),@(6,7),@(8,8),@(7,7),@(7,8),@(7,7),@(8,8),@(8,8),@(7,8),@(9,9),@(7,8),@(9,9),@(8,8),@(8,8),@(7,7),@(8,8),@(9,9),@(6,6),@(7,7),@(7,7),@(8,8),@(7,8),@(7,7),@(7,7),@(7,8),@(8,8),@(7,7),@(8,8),@(6,7),@(6,6),@(7,7),@(8,8),@(6,7),@(5,6),@(6,7),@(12,12),@(6,7),@(6,6),@(7,7),@(7,7),@(5,6),@(7,7),@(6,7),@(7,7),@(7,7),@(8,8),@(7,8),@(10,10),@(8,8),@(1,3),@(1,5),@(1,4),@(1,3),@(1,6),@(1,4),@(1,3),@(1,2))
for ($r = 2; $r -le 67; $r++) {
    $pair = $data[$r - 2]
    $ws.Cells.Item($r, 9).Value = $pair[0]
    $ws.Cells.Item($r, 10).Value = $pair[1]
}

Write-Output "Applied I0/IF columns"
